# Applies the diff:
#  1. Insert a new worksheet "Player Info" before "ODI Batting" with player
#     metadata (ID, NAME, BATTING_HAND, BOWL_STYLE).
#  2. On "ODI Batting": rename header D1 from MATCH_CARD_LINK -> MATCH_CODE
#     and change D2 from the full scorecard URL to just the match code "4284".

$wb = $excel.ActiveWorkbook

# --- 1. new "Player Info" sheet, inserted before the active ("ODI Batting") sheet ---
$info = $wb.Worksheets.Add()
$info.Name = "Player Info"

# --- existing sheet we will keep (becomes the 2nd tab); fetch AFTER the Add()
# call above so the reference resolves against the post-insert sheet layout ---
$odi = $wb.Worksheets.Item("ODI Batting")

$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

# Match the bold/bordered header style already used by the "ODI Batting"
# sheet's header row (copy formats only, values are left untouched).
$odi.Range("A1").Copy()
$info.Range("A1:D1").PasteSpecial(-4122)

# Keep the ID as text ("4832"), not a number, matching the source workbook.
$info.Range("A2").NumberFormat = "@"
$info.Range("A2").Value = "4832"
$info.Range("A2").Style = "Normal"

$info.Range("B2").Value = "Benjamin Thomas Foakes"
$info.Range("C2").Value = "Right Handed"
$info.Range("D2").Value = "Does Not Bowl | Unknown"

# --- 2. update "ODI Batting" sheet D column ---
$odi.Range("D1").Value = "MATCH_CODE"

# Keep the match code as text ("4284"), not a number, matching the source workbook.
$odi.Range("D2").NumberFormat = "@"
$odi.Range("D2").Value = "4284"
$odi.Range("D2").Style = "Normal"
